# "Visuals: Export to PNG, rename, reference in text"
#
# The first cluster of shapes on slide 15 (the picture/flowchart/connector/
# textbox group that makes up one "card" of the burning-questions visuals)
# was shifted 210004 EMU (~16.54pt) to the right, while keeping every
# shape's vertical position untouched. This corresponds to shapes 1-27 in
# that slide's z-order.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)

# Horizontal shift to apply, in EMU (English Metric Units; 914400 EMU = 1 in).
$deltaEmu = 210004

# Shape.Left/.Top round-trip through a single-precision float in this COM
# host, and converting points -> EMU on save truncates rather than rounds.
# Adding a small epsilon (well inside the safe margin between integer EMU
# boundaries) keeps the saved value exactly on the intended EMU target
# instead of landing 1 EMU short.
$eps = 0.000035

$shapeIndexes = @(1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27)

foreach ($i in $shapeIndexes) {
    $sh = $s.Shapes.Item($i)
    $oldEmu = [math]::Round($sh.Left * 12700)
    $newEmu = $oldEmu + $deltaEmu
    $sh.Left = ($newEmu / 12700.0) + $eps
}
